$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.768.97"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "1.864.79"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.039"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.18"
$ws.Range("E5").Value = "  +3.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.035"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  +3.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3802"
$ws.Range("E8").Value = "  +3.66%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07466"
$ws.Range("E9").Value = "  +3.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8851"
$ws.Range("E10").Value = "  +2.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.71"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "1.874.90"
$ws.Range("E12").Value = "  -15.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.563"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.767"
$ws.Range("E14").Value = "  +2.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07235"
$ws.Range("E15").Value = "  +4.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.80"
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.042"
$ws.Range("E17").Value = "  +2.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009161"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.035"
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.57"
$ws.Range("E20").Value = "  +2.67%  "
$ws.Range("D21").Value = "27.777.81"
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.997"
$ws.Range("E24").Value = "  +7.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.17"
$ws.Range("E25").Value = "  +3.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.88"
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.331"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("E28").Value = "  +4.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "117.99"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09067"
$ws.Range("E30").Value = "  +1.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7794"
$ws.Range("E31").Value = "  +4.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.097"
$ws.Range("E32").Value = "  +10.38%  "
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.574"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("E35").Value = "  +3.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.151"
$ws.Range("E36").Value = "  +2.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01994"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05350"
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.855"
$ws.Range("E39").Value = "  +3.64%  "
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1693"
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.905"
$ws.Range("E42").Value = "  +6.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.687"
$ws.Range("E43").Value = "  +4.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.85"
$ws.Range("E44").Value = "  +3.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.73"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.723"
$ws.Range("E46").Value = "  +4.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4712"
$ws.Range("E47").Value = "  +3.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06474"
$ws.Range("E48").Value = "  +4.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.923"
$ws.Range("E49").Value = "  +4.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.92"
$ws.Range("E50").Value = "  +2.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.60"
$ws.Range("E51").Value = "  +2.69%  "
